$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("G2").Value = 0.158418
$ws.Range("H2").Value = 0.475254
$ws.Range("I2").Value = 0.05652797120826585
$ws.Range("J2").Value = 0.05652797120826585
$ws.Range("M2").Value = 1.343359
$ws.Range("N2").Value = 4.030077
$ws.Range("O2").Value = 0.736296379391111
$ws.Range("P2").Value = 0.7362963793911109
$ws.Range("Q2").Value = 0.212812246062
$ws.Range("R2").Value = 1.915310214558
$ws.Range("S2").Value = 0.04162134053497111
$ws.Range("T2").Value = 0.0416213405349711

$ws.Range("G3").Value = 0.158418
$ws.Range("H3").Value = 0.475254
$ws.Range("I3").Value = 0.05652797120826585
$ws.Range("J3").Value = 0.05652797120826585
$ws.Range("O3").Value = 0.1764523396969075
$ws.Range("P3").Value = 0.1764523396969075
$ws.Range("Q3").Value = 0.05100014041199999
$ws.Range("R3").Value = 0.459001263708
$ws.Range("S3").Value = 0.00997449277801793
$ws.Range("T3").Value = 0.00997449277801793

$ws.Range("G4").Value = 0.158418
$ws.Range("H4").Value = 0.475254
$ws.Range("I4").Value = 0.05652797120826585
$ws.Range("J4").Value = 0.05652797120826585
$ws.Range("O4").Value = 0.08725128091198156
$ws.Range("P4").Value = 0.08725128091198156
$ws.Range("Q4").Value = 0.02521829739
$ws.Range("R4").Value = 0.22696467651
$ws.Range("S4").Value = 0.004932137895276809
$ws.Range("T4").Value = 0.004932137895276809

$ws.Range("I5").Value = 0.7628354881578912
$ws.Range("J5").Value = 0.7628354881578912
$ws.Range("M5").Value = 1.343359
$ws.Range("N5").Value = 4.030077
$ws.Range("O5").Value = 0.736296379391111
$ws.Range("P5").Value = 0.7362963793911109
$ws.Range("Q5").Value = 2.871865558602333
$ws.Range("R5").Value = 25.846790027421
$ws.Range("S5").Value = 0.561673008001706
$ws.Range("T5").Value = 0.5616730080017059

$ws.Range("I6").Value = 0.7628354881578912
$ws.Range("J6").Value = 0.7628354881578912
$ws.Range("O6").Value = 0.1764523396969075
$ws.Range("P6").Value = 0.1764523396969075
$ws.Range("S6").Value = 0.1346041066892924
$ws.Range("T6").Value = 0.1346041066892924

$ws.Range("I7").Value = 0.7628354881578912
$ws.Range("J7").Value = 0.7628354881578912
$ws.Range("O7").Value = 0.08725128091198156
$ws.Range("P7").Value = 0.08725128091198156
$ws.Range("S7").Value = 0.06655837346689274
$ws.Range("T7").Value = 0.06655837346689274

$ws.Range("G8").Value = 0.5062286666666667
$ws.Range("I8").Value = 0.180636540633843
$ws.Range("J8").Value = 0.180636540633843
$ws.Range("M8").Value = 1.343359
$ws.Range("N8").Value = 4.030077
$ws.Range("O8").Value = 0.736296379391111
$ws.Range("P8").Value = 0.7362963793911109
$ws.Range("Q8").Value = 0.6800468354246667
$ws.Range("R8").Value = 6.120421518822001
$ws.Range("S8").Value = 0.1330020308544339
$ws.Range("T8").Value = 0.1330020308544339

$ws.Range("G9").Value = 0.5062286666666667
$ws.Range("I9").Value = 0.180636540633843
$ws.Range("J9").Value = 0.180636540633843
$ws.Range("O9").Value = 0.1764523396969075
$ws.Range("P9").Value = 0.1764523396969075
$ws.Range("Q9").Value = 0.1629722195746667
$ws.Range("S9").Value = 0.0318737402295971
$ws.Range("T9").Value = 0.0318737402295971

$ws.Range("G10").Value = 0.5062286666666667
$ws.Range("I10").Value = 0.180636540633843
$ws.Range("J10").Value = 0.180636540633843
$ws.Range("O10").Value = 0.08725128091198156
$ws.Range("P10").Value = 0.08725128091198156
$ws.Range("Q10").Value = 0.08058569773222223
$ws.Range("R10").Value = 0.72527127959
$ws.Range("S10").Value = 0.01576076954981201
$ws.Range("T10").Value = 0.01576076954981201
